$wb = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item(1)

# Create two more copies of Sheet1, appended at the end, producing
# "Sheet1 (2)" and "Sheet1 (3)"
$ws1.Copy($null, $ws1)
$ws1.Copy($null, $wb.Worksheets.Item($wb.Worksheets.Count))

$ws2 = $wb.Worksheets.Item("Sheet1 (2)")
$ws3 = $wb.Worksheets.Item("Sheet1 (3)")

# ---- Edit "Sheet1 (2)" ----
$ws2.Activate()
$ws2.Range("C12").Value = "Blah"
$ws2.Range("E17").ClearContents()
$ws2.Range("B19").Value = "Ough"
[void]$ws2.Range("B20").Select()

# ---- Edit "Sheet1 (3)" ----
$ws3.Activate()
$ws3.Rows.Item(16).Insert()
$ws3.Range("B16").Value = "IM a new row"
$ws3.Range("B20").ClearContents()
$ws3.Range("D25").ClearContents()
[void]$ws3.Range("D26").Select()

# Restore Sheet1 selection/activation state and make "Sheet1 (3)" the
# active tab (matches activeTab="2" in the workbook views)
[void]$ws1.Range("D10").Select()
$ws3.Activate()
